$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new data rows (4 and 5) that replicate the existing pattern
# Row 4 mirrors row 3 ("ddaddress" pattern)
$ws.Range("A4").Value = "ddaddress"
$ws.Range("B4").Value = "dd city"
$ws.Range("C4").Value = " "
$ws.Range("D4").Value = "ddirst"
$ws.Range("E4").Value = "{{address}}"
$ws.Range("F4").Value = " "
$ws.Range("G4").Value = " "
$ws.Range("H4").Value = "{{address}}"

# Row 5 mirrors row 2 ("address" pattern)
$ws.Range("A5").Value = "address"
$ws.Range("B5").Value = "H C M City "
$ws.Range("C5").Value = " "
$ws.Range("D5").Value = "firsst"
$ws.Range("E5").Value = "last"
$ws.Range("F5").Value = " "
# "222" must land as text (matching the sibling cell on row 2), not a number,
# so enter it as a text formula then collapse it down to a literal value.
$g5 = $ws.Range("G5")
$g5.Formula = '="222"'
$g5.Copy()
$g5.PasteSpecial(-4163)
$ws.Range("H5").Value = " "
